# Update Sheets via scheduled runner: refresh market price / profit figures
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 8180.2144
$ws.Range("J2").Value = 16996.166
$ws.Range("L2").Value = 16996.166
$ws.Range("N2").Value = -17222.166
$ws.Range("H9").Value = 1845313.2
$ws.Range("I9").Value = 934.375
$ws.Range("J9").Value = 4796319.5
$ws.Range("K9").Value = 934.375
$ws.Range("L9").Value = 4796319.5
$ws.Range("M9").Value = -765.375
$ws.Range("N9").Value = -4796657.5
$ws.Range("H15").Value = 1515.25
$ws.Range("I15").Value = 1515.25
$ws.Range("K15").Value = 4545.75
$ws.Range("M15").Value = -4376.75
$ws.Range("H29").Value = 8857.143
$ws.Range("I29").Value = 664.3333
$ws.Range("J29").Value = 15001.75
$ws.Range("K29").Value = 1992.9999
$ws.Range("L29").Value = 45005.25
$ws.Range("M29").Value = -1711.9999
$ws.Range("N29").Value = -45567.25
$ws.Range("H38").Value = 559.9286
$ws.Range("I38").Value = 218.38461
$ws.Range("K38").Value = 655.15383
$ws.Range("M38").Value = -283.15383
$ws.Range("H80").Value = 833.2222
$ws.Range("I80").Value = 502.2
$ws.Range("J80").Value = 1247
$ws.Range("K80").Value = 1506.6
$ws.Range("L80").Value = 3741
$ws.Range("M80").Value = -508.5999999999999
$ws.Range("N80").Value = -5737
$ws.Range("H83").Value = 833.2222
$ws.Range("I83").Value = 502.2
$ws.Range("J83").Value = 1247
$ws.Range("K83").Value = 4519.8
$ws.Range("L83").Value = 11223
$ws.Range("M83").Value = 472.1999999999998
$ws.Range("N83").Value = -21207
$ws.Range("H132").Value = 1194.5122
$ws.Range("I132").Value = 1147.9744
$ws.Range("K132").Value = 3443.9232
$ws.Range("M132").Value = -913.9232000000002
$ws.Range("H137").Value = 4242619.5
$ws.Range("I137").Value = 6581335.5
$ws.Range("J137").Value = 10657.857
$ws.Range("K137").Value = 19744006.5
$ws.Range("L137").Value = 31973.571
$ws.Range("M137").Value = -19741456.5
$ws.Range("N137").Value = -37073.571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4882785.5
$ws.Range("I32").Value = 6455808
$ws.Range("K32").Value = 6455808
$ws.Range("M32").Value = -6455521
$ws.Range("H45").Value = 2273
$ws.Range("I45").Value = 2312.0715
$ws.Range("J45").Value = 1999.5
$ws.Range("K45").Value = 2312.0715
$ws.Range("L45").Value = 1999.5
$ws.Range("M45").Value = -1935.0715
$ws.Range("N45").Value = -2753.5
$ws.Range("H61").Value = 4604
$ws.Range("I61").Value = 3219.65
$ws.Range("K61").Value = 3219.65
$ws.Range("M61").Value = -3007.65
$ws.Range("H122").Value = 3937.182
$ws.Range("I122").Value = 3701.3333
$ws.Range("K122").Value = 11103.9999
$ws.Range("M122").Value = -8653.999899999999
$ws.Range("H132").Value = 11559
$ws.Range("I132").Value = 6124.75
$ws.Range("K132").Value = 18374.25
$ws.Range("M132").Value = -15844.25
$ws.Range("H136").Value = 4604
$ws.Range("I136").Value = 3219.65
$ws.Range("K136").Value = 9658.950000000001
$ws.Range("M136").Value = -7108.950000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 3266.6667
$ws.Range("I36").Value = 3266.6667
$ws.Range("K36").Value = 3266.6667
$ws.Range("M36").Value = -2732.6667
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("N39").ClearContents()
$ws.Range("H94").Value = 672.93335
$ws.Range("I94").Value = 554.5
$ws.Range("K94").Value = 554.5
$ws.Range("M94").Value = -103.5
$ws.Range("H107").Value = 2256.186
$ws.Range("I107").Value = 1380.7297
$ws.Range("J107").Value = 7654.8335
$ws.Range("K107").Value = 1380.7297
$ws.Range("L107").Value = 7654.8335
$ws.Range("M107").Value = 539.2702999999999
$ws.Range("N107").Value = -11494.8335
$ws.Range("H140").Value = 76922.62
$ws.Range("J140").Value = 76922.62
$ws.Range("L140").Value = 76922.62
$ws.Range("N140").Value = -87282.62

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4522.2
$ws.Range("I16").Value = 4203.6665
$ws.Range("J16").Value = 5000
$ws.Range("K16").Value = 4203.6665
$ws.Range("L16").Value = 5000
$ws.Range("M16").Value = -3916.6665
$ws.Range("N16").Value = -5574
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("H58").Value = 7430.625
$ws.Range("I58").Value = 5708.4443
$ws.Range("K58").Value = 5708.4443
$ws.Range("M58").Value = -5505.4443
$ws.Range("H113").Value = 4522.2
$ws.Range("I113").Value = 4203.6665
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 4203.6665
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = -2033.6665
$ws.Range("N113").Value = -9340
$ws.Range("H132").Value = 49470.11
$ws.Range("I132").Value = 3919.8572
$ws.Range("J132").Value = 98524.234
$ws.Range("K132").Value = 11759.5716
$ws.Range("L132").Value = 295572.702
$ws.Range("M132").Value = -9229.571599999999
$ws.Range("N132").Value = -300632.702
$ws.Range("H136").Value = 7430.625
$ws.Range("I136").Value = 5708.4443
$ws.Range("K136").Value = 17125.3329
$ws.Range("M136").Value = -14575.3329

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 5805
$ws.Range("I56").Value = 5805
$ws.Range("K56").Value = 5805
$ws.Range("M56").Value = -5275
$ws.Range("H87").Value = 1240
$ws.Range("I87").Value = 1240
$ws.Range("K87").Value = 3720
$ws.Range("M87").Value = -2472
$ws.Range("H90").Value = 1240
$ws.Range("I90").Value = 1240
$ws.Range("K90").Value = 11160
$ws.Range("M90").Value = -4920
$ws.Range("H114").Value = 563
$ws.Range("J114").Value = 1000
$ws.Range("L114").Value = 3000
$ws.Range("N114").Value = -9508
$ws.Range("H136").Value = 3039.6667
$ws.Range("I136").Value = 2684.182
$ws.Range("J136").Value = 6950
$ws.Range("K136").Value = 8052.545999999999
$ws.Range("L136").Value = 20850
$ws.Range("M136").Value = -2952.545999999999
$ws.Range("N136").Value = -31050

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7881.05
$ws.Range("I70").Value = 8274.143
$ws.Range("J70").Value = 6963.8335
$ws.Range("K70").Value = 8274.143
$ws.Range("L70").Value = 6963.8335
$ws.Range("M70").Value = -8004.143
$ws.Range("N70").Value = -7503.8335
$ws.Range("H73").Value = 7881.05
$ws.Range("I73").Value = 8274.143
$ws.Range("J73").Value = 6963.8335
$ws.Range("K73").Value = 8274.143
$ws.Range("L73").Value = 6963.8335
$ws.Range("M73").Value = -7338.143
$ws.Range("N73").Value = -8835.833500000001
$ws.Range("H122").Value = 11974.571
$ws.Range("I122").Value = 12620.333
$ws.Range("J122").Value = 8100
$ws.Range("K122").Value = 37860.999
$ws.Range("L122").Value = 24300
$ws.Range("M122").Value = -35410.999
$ws.Range("N122").Value = -29200
$ws.Range("H132").Value = 11038.889
$ws.Range("I132").Value = 1680.5
$ws.Range("J132").Value = 13712.714
$ws.Range("K132").Value = 5041.5
$ws.Range("L132").Value = 41138.142
$ws.Range("M132").Value = -2511.5
$ws.Range("N132").Value = -46198.142

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1898.3334
$ws.Range("I16").Value = 800.3333
$ws.Range("J16").Value = 2996.3333
$ws.Range("K16").Value = 800.3333
$ws.Range("L16").Value = 2996.3333
$ws.Range("M16").Value = -630.3333
$ws.Range("N16").Value = -3336.3333
$ws.Range("H20").Value = 25000
$ws.Range("J20").Value = 25000
$ws.Range("L20").Value = 25000
$ws.Range("N20").Value = -25452
$ws.Range("H22").Value = 2325.6553
$ws.Range("J22").Value = 2793.6667
$ws.Range("L22").Value = 2793.6667
$ws.Range("N22").Value = -3383.6667
$ws.Range("H27").Value = 2325.6553
$ws.Range("J27").Value = 2793.6667
$ws.Range("L27").Value = 2793.6667
$ws.Range("N27").Value = -3007.6667
$ws.Range("H40").Value = 100010360
$ws.Range("I40").Value = 100010360
$ws.Range("K40").Value = 100010360
$ws.Range("M40").Value = -100010224
$ws.Range("H55").Value = 637.9524
$ws.Range("I55").Value = 409.07693
$ws.Range("J55").Value = 1009.875
$ws.Range("K55").Value = 409.07693
$ws.Range("L55").Value = 1009.875
$ws.Range("M55").Value = -236.07693
$ws.Range("N55").Value = -1355.875
$ws.Range("H100").Value = 9263063
$ws.Range("I100").Value = 13892053
$ws.Range("K100").Value = 13892053
$ws.Range("M100").Value = -13891512
$ws.Range("H122").Value = 35718070
$ws.Range("I122").Value = 50003510
$ws.Range("J122").Value = 4472.25
$ws.Range("K122").Value = 150010530
$ws.Range("L122").Value = 13416.75
$ws.Range("M122").Value = -150008080
$ws.Range("N122").Value = -18316.75
$ws.Range("H136").Value = 5994.8
$ws.Range("I136").Value = 3211
$ws.Range("K136").Value = 9633
$ws.Range("M136").Value = -7083

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 16618.875
$ws.Range("J74").Value = 16618.875
$ws.Range("L74").Value = 16618.875
$ws.Range("N74").Value = -18490.875
$ws.Range("H77").Value = 16618.875
$ws.Range("J77").Value = 16618.875
$ws.Range("L77").Value = 49856.625
$ws.Range("N77").Value = -59216.625
$ws.Range("H132").Value = 4642.7617
$ws.Range("I132").Value = 2549.8667
$ws.Range("K132").Value = 7649.6001
$ws.Range("M132").Value = -5119.6001
$ws.Range("H135").Value = 88166.5
$ws.Range("J135").Value = 102999.664
$ws.Range("L135").Value = 102999.664
$ws.Range("N135").Value = -113139.664
$ws.Range("H136").Value = 5429.811
$ws.Range("I136").Value = 4827.2173
$ws.Range("J136").Value = 6419.7856
$ws.Range("K136").Value = 14481.6519
$ws.Range("L136").Value = 19259.3568
$ws.Range("M136").Value = -11931.6519
$ws.Range("N136").Value = -24359.3568

Write-Host "Updated 256 cells across 8 sheets"
